$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Devices")

# Refresh the UPDATED_TS column (K) for a handful of devices as part of a
# test iteration run.
$ws.Range("K2").Value = "2025-07-10T20:32:30.372549546+02:00[Europe/Amsterdam]"
$ws.Range("K3").Value = "2025-07-10T12:16:59.450154591+02:00[Europe/Amsterdam]"
$ws.Range("K16").Value = "2025-07-10T12:16:59.729918525+02:00[Europe/Amsterdam]"
$ws.Range("K26").Value = "2025-07-10T20:32:27.944336134+02:00[Europe/Amsterdam]"
$ws.Range("K37").Value = "2025-07-10T20:32:29.898247690+02:00[Europe/Amsterdam]"

# Row 49 (SMART_LIGHT / SL001) was re-registered: brand/model cleared, the
# supported actions list lost "setMode", and UPDATED_TS bumped.
$ws.Range("D49").Value = ""
$ws.Range("E49").Value = ""
$ws.Range("I49").Value = "on, off, status"
$ws.Range("K49").Value = "2025-07-09T22:58:25.093936180+02:00[Europe/Amsterdam]"
